# This script applies a row-permutation update to the "Fruta, Agrícola del
# Norte S.A. de Arica - Tuna" weekly price sheet. The diff shows that the
# data rows (2-12 and 15-18) were reshuffled: each row's Fecha/Calidad/
# Volumen/Precio.../Unidad/Origen/Precio-Kg/Kg-unidad block now appears at a
# different row position, while the identifying columns (A,B,C,E-K) stay the
# same for every row. We therefore rewrite, for each target row, the full
# set of values that (per the diff) now belongs there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2 now holds the data that used to be in row 11
$ws.Range("D2").Value = 44679
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 29000
$ws.Range("O2").Value = 30000
$ws.Range("P2").Value = 29500
$ws.Range("Q2").Value = "$/caja 20 kilos"
$ws.Range("R2").Value = "Región de Coquimbo"
$ws.Range("S2").Value = 1475
$ws.Range("T2").Value = 20

# Row 3 now holds the data that used to be in row 12
$ws.Range("D3").Value = 44679
$ws.Range("L3").Value = "Tercera"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 24000
$ws.Range("O3").Value = 25000
$ws.Range("P3").Value = 24500
$ws.Range("Q3").Value = "$/caja 20 kilos"
$ws.Range("R3").Value = "Región de Coquimbo"
$ws.Range("S3").Value = 1225
$ws.Range("T3").Value = 20

# Row 4 now holds the data that used to be in row 2
$ws.Range("D4").Value = 45028
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 21000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 21500
$ws.Range("Q4").Value = "$/caja 20 kilos"
$ws.Range("R4").Value = "Región de Coquimbo"
$ws.Range("S4").Value = 1075
$ws.Range("T4").Value = 20

# Row 5 now holds the data that used to be in row 6
$ws.Range("D5").Value = 44650
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 160
$ws.Range("N5").Value = 31000
$ws.Range("O5").Value = 32000
$ws.Range("P5").Value = 31500
$ws.Range("Q5").Value = "$/caja 20 kilos"
$ws.Range("R5").Value = "Región de Coquimbo"
$ws.Range("S5").Value = 1575
$ws.Range("T5").Value = 20

# Row 6 now holds the data that used to be in row 7
$ws.Range("D6").Value = 44650
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 29000
$ws.Range("O6").Value = 30000
$ws.Range("P6").Value = 29500
$ws.Range("Q6").Value = "$/caja 20 kilos"
$ws.Range("R6").Value = "Región de Coquimbo"
$ws.Range("S6").Value = 1475
$ws.Range("T6").Value = 20

# Row 7 now holds the data that used to be in row 8
$ws.Range("D7").Value = 44979
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 250
$ws.Range("N7").Value = 29000
$ws.Range("O7").Value = 30000
$ws.Range("P7").Value = 29500
$ws.Range("Q7").Value = "$/caja 20 kilos"
$ws.Range("R7").Value = "Región de Coquimbo"
$ws.Range("S7").Value = 1475
$ws.Range("T7").Value = 20

# Row 8 now holds the data that used to be in row 16
$ws.Range("D8").Value = 44671
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 29000
$ws.Range("O8").Value = 30000
$ws.Range("P8").Value = 29500
$ws.Range("Q8").Value = "$/caja 20 kilos"
$ws.Range("R8").Value = "Región de Coquimbo"
$ws.Range("S8").Value = 1475
$ws.Range("T8").Value = 20

# Row 9 now holds the data that used to be in row 15
$ws.Range("D9").Value = 45007
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 160
$ws.Range("N9").Value = 27000
$ws.Range("O9").Value = 28000
$ws.Range("P9").Value = 27500
$ws.Range("Q9").Value = "$/caja 20 kilos"
$ws.Range("R9").Value = "Región de Coquimbo"
$ws.Range("S9").Value = 1375
$ws.Range("T9").Value = 20

# Row 10 now holds the data that used to be in row 5
$ws.Range("D10").Value = 45021
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 250
$ws.Range("N10").Value = 22000
$ws.Range("O10").Value = 23000
$ws.Range("P10").Value = 22500
$ws.Range("Q10").Value = "$/caja 20 kilos"
$ws.Range("R10").Value = "Región de Coquimbo"
$ws.Range("S10").Value = 1125
$ws.Range("T10").Value = 20

# Row 11 now holds the data that used to be in row 17
$ws.Range("D11").Value = 44972
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 140
$ws.Range("N11").Value = 27000
$ws.Range("O11").Value = 28000
$ws.Range("P11").Value = 27429
$ws.Range("Q11").Value = "$/caja 18 kilos"
$ws.Range("R11").Value = "Región Metropolitana"
$ws.Range("S11").Value = 1524
$ws.Range("T11").Value = 18

# Row 12 now holds the data that used to be in row 3
$ws.Range("D12").Value = 44993
$ws.Range("L12").Value = "Segunda"
$ws.Range("M12").Value = 130
$ws.Range("N12").Value = 25000
$ws.Range("O12").Value = 26000
$ws.Range("P12").Value = 25462
$ws.Range("Q12").Value = "$/caja 20 kilos"
$ws.Range("R12").Value = "Región de Coquimbo"
$ws.Range("S12").Value = 1273
$ws.Range("T12").Value = 20

# Row 15 now holds the data that used to be in row 18
$ws.Range("D15").Value = 44643
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 160
$ws.Range("N15").Value = 28000
$ws.Range("O15").Value = 30000
$ws.Range("P15").Value = 29000
$ws.Range("Q15").Value = "$/caja 20 kilos"
$ws.Range("R15").Value = "Región de Coquimbo"
$ws.Range("S15").Value = 1450
$ws.Range("T15").Value = 20

# Row 16 now holds the data that used to be in row 4
$ws.Range("D16").Value = 45014
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 24000
$ws.Range("O16").Value = 25000
$ws.Range("P16").Value = 24500
$ws.Range("Q16").Value = "$/caja 20 kilos"
$ws.Range("R16").Value = "Región de Coquimbo"
$ws.Range("S16").Value = 1225
$ws.Range("T16").Value = 20

# Row 17 now holds the data that used to be in row 9
$ws.Range("D17").Value = 44636
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 29000
$ws.Range("O17").Value = 30000
$ws.Range("P17").Value = 29500
$ws.Range("Q17").Value = "$/caja 20 kilos"
$ws.Range("R17").Value = "Región de Coquimbo"
$ws.Range("S17").Value = 1475
$ws.Range("T17").Value = 20

# Row 18 now holds the data that used to be in row 10
$ws.Range("D18").Value = 44664
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 150
$ws.Range("N18").Value = 29000
$ws.Range("O18").Value = 30000
$ws.Range("P18").Value = 29500
$ws.Range("Q18").Value = "$/caja 18 kilos"
$ws.Range("R18").Value = "Región de Coquimbo"
$ws.Range("S18").Value = 1639
$ws.Range("T18").Value = 18

Write-Host "Applied weekly Tuna price update (row permutation)."
